$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3124
$ws1.Range("F4").Value = 1081
$ws1.Range("F5").Value = 78
$ws1.Range("F7").Value = 270
$ws1.Range("F8").Value = 32
$ws1.Range("F9").Value = 1112
$ws1.Range("F10").Value = 15522
$ws1.Range("F11").Value = 232
$ws1.Range("F12").Value = 164
$ws1.Range("F13").Value = 1018
$ws1.Range("F14").Value = 6129
$ws1.Range("F15").Value = 620
$ws1.Range("F16").Value = 104
$ws1.Range("F17").Value = 63
$ws1.Range("F18").Value = 4
$ws1.Range("F19").Value = 107
$ws1.Range("F22").Value = 218
$ws1.Range("F26").Value = 208
$ws1.Range("F27").Value = 856
$ws1.Range("F28").Value = 20
$ws1.Range("F29").Value = 4988
$ws1.Range("F30").Value = 415
$ws1.Range("F31").Value = 10993
$ws1.Range("F32").Value = 1225
$ws1.Range("F35").Value = 154
$ws1.Range("F36").Value = 3791
$ws1.Range("F37").Value = 262

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3124
$ws4.Range("F5").Value = 1081
$ws4.Range("F6").Value = 78
$ws4.Range("F8").Value = 270
$ws4.Range("F9").Value = 32
$ws4.Range("F10").Value = 1112
$ws4.Range("F11").Value = 15523
$ws4.Range("F12").Value = 232
$ws4.Range("F13").Value = 164
$ws4.Range("F14").Value = 1018
$ws4.Range("F15").Value = 6129
$ws4.Range("F16").Value = 620
$ws4.Range("F17").Value = 104
$ws4.Range("F18").Value = 63
$ws4.Range("F19").Value = 4
$ws4.Range("F20").Value = 107
$ws4.Range("F23").Value = 218
$ws4.Range("F27").Value = 208
$ws4.Range("F28").Value = 856
$ws4.Range("F29").Value = 20
$ws4.Range("F30").Value = 4988
$ws4.Range("F31").Value = 415
$ws4.Range("F33").Value = 10993
$ws4.Range("F34").Value = 1225
$ws4.Range("F37").Value = 154
$ws4.Range("F38").Value = 3791
$ws4.Range("F39").Value = 262
